$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the original row values (as displayed text) before anything is
# overwritten. Columns F1 and I1 were genuinely empty cells in the source
# sheet (string-typed cell with no <v>), so they are hard-coded as "" rather
# than re-read (the COM layer mis-resolves a missing <v> to shared-string 0).
$origB = $ws.Range("B1").Text
$origC = $ws.Range("C1").Text
$origD = $ws.Range("D1").Text
$origE = $ws.Range("E1").Text
$origF = ""
$origG = $ws.Range("G1").Text
$origH = $ws.Range("H1").Text
$origI = ""
$origJ = $ws.Range("J1").Text
$origK = $ws.Range("K1").Text
$origL = $ws.Range("L1").Text
$origM = $ws.Range("M1").Text

# Force the destination cells to stay text so re-assigning values such as
# "2018-10-14" or "6" does not get reinterpreted as a date/number.
$ws.Range("A1:O1").NumberFormat = "@"

# Shift the existing row (B1:M1) two columns to the right, into D1:O1, to
# make room for the new player-identity columns.
$ws.Range("D1").Value = $origB
$ws.Range("E1").Value = $origC
$ws.Range("F1").Value = $origD
$ws.Range("G1").Value = $origE
$ws.Range("H1").Value = $origF
$ws.Range("I1").Value = $origG
$ws.Range("J1").Value = $origH
$ws.Range("K1").Value = $origI
$ws.Range("L1").Value = $origJ
$ws.Range("M1").Value = $origK
$ws.Range("N1").Value = $origL
$ws.Range("O1").Value = $origM

# New leading player-identity columns.
$ws.Range("A1").Value = "Kelly"
$ws.Range("B1").Value = "Chad"
$ws.Range("C1").Value = "QB"

# Drop the temporary text formatting so styles.xml is left untouched.
$ws.Range("A1:O1").ClearFormats()

# New trailing numeric column.
$ws.Range("P1").Value = 0
